$wb = $excel.ActiveWorkbook

# Sheet "展览" (sheet1) - F column ("想去人数") updates
$ws1 = $wb.Worksheets.Item(1)
$ws1.Cells.Item(2, 6).Value = 29
$ws1.Cells.Item(3, 6).Value = 222
$ws1.Cells.Item(4, 6).Value = 4763
$ws1.Cells.Item(7, 6).Value = 110
$ws1.Cells.Item(8, 6).Value = 101
$ws1.Cells.Item(10, 6).Value = 736
$ws1.Cells.Item(11, 6).Value = 216
$ws1.Cells.Item(12, 6).Value = 1129
$ws1.Cells.Item(13, 6).Value = 100
$ws1.Cells.Item(14, 6).Value = 259
$ws1.Cells.Item(15, 6).Value = 175
$ws1.Cells.Item(17, 6).Value = 141
$ws1.Cells.Item(19, 6).Value = 3867
$ws1.Cells.Item(20, 6).Value = 6196
$ws1.Cells.Item(23, 6).Value = 84
$ws1.Cells.Item(26, 6).Value = 3947
$ws1.Cells.Item(27, 6).Value = 395
$ws1.Cells.Item(28, 6).Value = 36
$ws1.Cells.Item(29, 6).Value = 2542
$ws1.Cells.Item(32, 6).Value = 0
$ws1.Cells.Item(33, 6).Value = 269
$ws1.Cells.Item(34, 6).Value = 295
$ws1.Cells.Item(35, 6).Value = 364
$ws1.Cells.Item(36, 6).Value = 161
$ws1.Cells.Item(37, 6).Value = 1555
$ws1.Cells.Item(39, 6).Value = 42
$ws1.Cells.Item(42, 6).Value = 484
$ws1.Cells.Item(43, 6).Value = 478
$ws1.Cells.Item(44, 6).Value = 73
$ws1.Cells.Item(45, 6).Value = 572

# Sheet "全部类型" (sheet4) - F column ("想去人数") updates
$ws4 = $wb.Worksheets.Item(4)
$ws4.Cells.Item(3, 6).Value = 222
$ws4.Cells.Item(4, 6).Value = 4763
$ws4.Cells.Item(6, 6).Value = 146
$ws4.Cells.Item(8, 6).Value = 106
$ws4.Cells.Item(9, 6).Value = 101
$ws4.Cells.Item(11, 6).Value = 736
$ws4.Cells.Item(12, 6).Value = 216
$ws4.Cells.Item(13, 6).Value = 1129
$ws4.Cells.Item(14, 6).Value = 100
$ws4.Cells.Item(16, 6).Value = 175
$ws4.Cells.Item(18, 6).Value = 141
$ws4.Cells.Item(20, 6).Value = 3867
$ws4.Cells.Item(21, 6).Value = 6196
$ws4.Cells.Item(22, 6).Value = 39
$ws4.Cells.Item(23, 6).Value = 38
$ws4.Cells.Item(24, 6).Value = 84
$ws4.Cells.Item(25, 6).Value = 0
$ws4.Cells.Item(26, 6).Value = 46
$ws4.Cells.Item(27, 6).Value = 3947
$ws4.Cells.Item(28, 6).Value = 395
$ws4.Cells.Item(29, 6).Value = 36
$ws4.Cells.Item(30, 6).Value = 2542
$ws4.Cells.Item(31, 6).Value = 568
$ws4.Cells.Item(33, 6).Value = 137
$ws4.Cells.Item(34, 6).Value = 269
$ws4.Cells.Item(35, 6).Value = 295
$ws4.Cells.Item(36, 6).Value = 364
$ws4.Cells.Item(37, 6).Value = 161
$ws4.Cells.Item(38, 6).Value = 1555
$ws4.Cells.Item(43, 6).Value = 484
$ws4.Cells.Item(46, 6).Value = 572

